$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the old "Hoja1" sheet, leaving "DatosCP" as the sole worksheet ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Delete() | Out-Null

# --- Update the "DatosCP" sheet: first data case (row 2) replaced with the new test case ---
$ws = $wb.Worksheets.Item("DatosCP")
$ws.Activate()

$ws.Range("B2").Value = "qweqweew"
$ws.Range("C2").Value = "ee51165"
$ws.Range("A2").Value = "CP001_loginInvalidEmail"
$ws.Range("D2").Value = "Invalid email address."

# The new test-case id (A2) carries its own style: underlined font, same text
# number format + border as the rest of the column.
$ws.Range("A2").Font.Underline = 2

$excel.DisplayAlerts = $true
